# Update trading results - Sat Sep 20 12:39:24 UTC 2025
# Append two new trading log rows (24 and 25) to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: TRADING_ATTEMPT
$ws.Cells.Item(24, 1).Value = "2025-09-20T12:39:21.948063"
$ws.Cells.Item(24, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(24, 3).Value = "ENA"
$ws.Cells.Item(24, 4).Value = "UNKNOWN"
$ws.Cells.Item(24, 5).Value = 0.6727811902747289
$ws.Cells.Item(24, 6).Value = ""
$ws.Cells.Item(24, 7).Value = ""
$ws.Cells.Item(24, 8).Value = ""
$ws.Cells.Item(24, 9).Value = ""
$ws.Cells.Item(24, 10).Value = ""
$ws.Cells.Item(24, 11).Value = "ATTEMPT"
$ws.Cells.Item(24, 12).Value = "Attempting trade 1/1"

# Row 25: POSITION_OPENED
$ws.Cells.Item(25, 1).Value = "2025-09-20T12:39:23.406345"
$ws.Cells.Item(25, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(25, 3).Value = "ENA"
$ws.Cells.Item(25, 4).Value = "UNKNOWN"
$ws.Cells.Item(25, 5).Value = 0.6727811902747289
$ws.Cells.Item(25, 6).Value = 1200
$ws.Cells.Item(25, 7).Value = 10
$ws.Cells.Item(25, 8).Value = 0.07577367435598509
$ws.Cells.Item(25, 9).Value = ""
$ws.Cells.Item(25, 10).Value = ""
$ws.Cells.Item(25, 11).Value = "SUCCESS"
$ws.Cells.Item(25, 12).Value = ""
